$d = $word.ActiveDocument

# Find the paragraph that contains the "Ver no Jupiter" text and the one with
# the "(c) 2020" copyright text, together with the blank paragraph that
# precedes them, and remove all three paragraphs (including their paragraph
# marks) while leaving the preceding "de materiais..." paragraph and the
# following blank paragraph untouched.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*de materiais e de equipamentos*") {
        $startPara = $p
    }
    if ($t -like "*2020*Contact*luizeleno*") {
        $endPara = $p
    }
}

$start = $startPara.Range.End
$end = $endPara.Range.End

$r = $d.Range($start, $end)
$r.Delete()
